$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "IP" header (H1) onto the two new
# header cells so they pick up the same style (bold, bordered, centered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Add new headers in I1 and J1
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Fill in I and J columns for data rows 2-24
# Column I ("I0") is always 1
# Column J ("IF") mirrors the value in column H for that row
for ($row = 2; $row -le 24; $row++) {
    $hValue = $ws.Cells.Item($row, 8).Value2
    $ws.Cells.Item($row, 9).Value2 = 1
    $ws.Cells.Item($row, 10).Value2 = $hValue
}
